$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '29.278.57'
$ws.Cells.Item(2, 5).Value = '  -0.33%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.863.15'
$ws.Cells.Item(3, 5).Value = '  -1.14%  '

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.004'
$ws.Cells.Item(4, 5).Value = '  +0.44%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '242.56'
$ws.Cells.Item(5, 5).Value = '  +0.05%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.7019'
$ws.Cells.Item(6, 5).Value = '  -1.70%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '1.004'
$ws.Cells.Item(7, 5).Value = '  +0.44%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.07786'

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.3102'
$ws.Cells.Item(9, 5).Value = '  -1.30%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '24.19'
$ws.Cells.Item(10, 5).Value = '  -4.65%  '

# Row 11
$ws.Cells.Item(11, 2).Value = 'WrappedEther'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(11, 4).Value = '2.418.05'
$ws.Cells.Item(11, 5).Value = '  +29.03%  '

# Row 12
$ws.Cells.Item(12, 2).Value = 'TRON'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.08020'
$ws.Cells.Item(12, 5).Value = '  -4.12%  '

# Row 13
$ws.Cells.Item(13, 2).Value = 'Polkadot'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '5.167'
$ws.Cells.Item(13, 5).Value = '  -1.69%  '

# Row 14
$ws.Cells.Item(14, 2).Value = 'Litecoin'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '92.79'
$ws.Cells.Item(14, 5).Value = '  +1.02%  '

# Row 15
$ws.Cells.Item(15, 2).Value = 'Polygon'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.6943'
$ws.Cells.Item(15, 5).Value = '  -4.01%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '6.320'
$ws.Cells.Item(16, 5).Value = '  +0.51%  '

# Row 17
$ws.Cells.Item(17, 2).Value = 'ShibaInu'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.000008236'
$ws.Cells.Item(17, 5).Value = '  -2.70%  '

# Row 18
$ws.Cells.Item(18, 2).Value = 'WrappedBTC'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(18, 4).Value = '28.580.50'
$ws.Cells.Item(18, 5).Value = '  -2.66%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '249.01'
$ws.Cells.Item(19, 5).Value = '  +2.96%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '2.122.49'
$ws.Cells.Item(20, 5).Value = '  +0.67%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '13.12'
$ws.Cells.Item(21, 5).Value = '  -1.15%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '1.001'
$ws.Cells.Item(22, 5).Value = '  +0.14%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '7.498'
$ws.Cells.Item(23, 5).Value = '  -4.01%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '1.004'
$ws.Cells.Item(24, 5).Value = '  +0.41%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  -2.61%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '8.963'
$ws.Cells.Item(26, 5).Value = '  -1.39%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '160.06'

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '18.57'
$ws.Cells.Item(28, 5).Value = '  -0.17%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.501'
$ws.Cells.Item(29, 5).Value = '  -0.46%  '

# Row 30
$ws.Cells.Item(30, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '4.260'
$ws.Cells.Item(30, 5).Value = '  -2.31%  '

# Row 31
$ws.Cells.Item(31, 2).Value = 'Filecoin'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '4.266'
$ws.Cells.Item(31, 5).Value = '  -3.83%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '1.221'
$ws.Cells.Item(32, 5).Value = '  +0.61%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.05235'
$ws.Cells.Item(33, 5).Value = '  -2.79%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.879'
$ws.Cells.Item(34, 5).Value = '  -3.95%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.7382'
$ws.Cells.Item(35, 5).Value = '  -2.09%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '1.157'
$ws.Cells.Item(36, 5).Value = '  -2.07%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '2.711'
$ws.Cells.Item(37, 5).Value = '  +0.60%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.01856'
$ws.Cells.Item(38, 5).Value = '  -1.56%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '1.255.63'
$ws.Cells.Item(39, 5).Value = '  -2.05%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '2.741'
$ws.Cells.Item(40, 5).Value = '  -0.09%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '6.174'
$ws.Cells.Item(41, 5).Value = '  -6.19%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.8951'
$ws.Cells.Item(42, 5).Value = '  +0.32%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'Quant'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '110.62'
$ws.Cells.Item(43, 5).Value = '  +0.07%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '71.32'
$ws.Cells.Item(44, 5).Value = '  -2.92%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '1.004'
$ws.Cells.Item(45, 5).Value = '  +0.40%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.00000000129'
$ws.Cells.Item(46, 5).Value = '  -0.80%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Cells.Item(47, 4).Value = '2.021.72'
$ws.Cells.Item(47, 5).Value = '  +0.88%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.5206'
$ws.Cells.Item(48, 5).Value = '  -0.22%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '1.786'
$ws.Cells.Item(49, 5).Value = '  -1.23%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  -2.01%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '1.003'
$ws.Cells.Item(51, 5).Value = '  +0.36%  '
